$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 864
$ws1.Cells.Item(3, 6).Value = 1439
$ws1.Cells.Item(4, 6).Value = 1090
$ws1.Cells.Item(5, 6).Value = 516
$ws1.Cells.Item(7, 6).Value = 662
$ws1.Cells.Item(8, 6).Value = 240
$ws1.Cells.Item(9, 6).Value = 10
$ws1.Cells.Item(10, 6).Value = 80
$ws1.Cells.Item(11, 6).Value = 215
$ws1.Cells.Item(12, 6).Value = 148
$ws1.Cells.Item(13, 6).Value = 1789
$ws1.Cells.Item(14, 6).Value = 427
$ws1.Cells.Item(16, 6).Value = 488
$ws1.Cells.Item(17, 6).Value = 255
$ws1.Cells.Item(19, 6).Value = 113
$ws1.Cells.Item(21, 6).Value = 659
$ws1.Cells.Item(22, 6).Value = 46
$ws1.Cells.Item(23, 6).Value = 239
$ws1.Cells.Item(24, 6).Value = 955
$ws1.Cells.Item(26, 6).Value = 1541
$ws1.Cells.Item(27, 6).Value = 285

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 121
$ws2.Cells.Item(5, 6).Value = 208

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 864
$ws4.Cells.Item(4, 6).Value = 1439
$ws4.Cells.Item(5, 6).Value = 1090
$ws4.Cells.Item(6, 6).Value = 121
$ws4.Cells.Item(8, 6).Value = 516
$ws4.Cells.Item(10, 6).Value = 662
$ws4.Cells.Item(12, 6).Value = 240
$ws4.Cells.Item(13, 6).Value = 10
$ws4.Cells.Item(14, 6).Value = 80
$ws4.Cells.Item(15, 6).Value = 215
$ws4.Cells.Item(16, 6).Value = 148
$ws4.Cells.Item(17, 6).Value = 1789
$ws4.Cells.Item(18, 6).Value = 208
$ws4.Cells.Item(19, 6).Value = 427
$ws4.Cells.Item(21, 6).Value = 488
$ws4.Cells.Item(22, 6).Value = 255
$ws4.Cells.Item(25, 6).Value = 113
$ws4.Cells.Item(30, 6).Value = 659
$ws4.Cells.Item(35, 6).Value = 46
$ws4.Cells.Item(36, 6).Value = 239
$ws4.Cells.Item(37, 6).Value = 955
$ws4.Cells.Item(39, 6).Value = 1541
$ws4.Cells.Item(40, 6).Value = 285
